# apa6_man.docx template update:
#  1. Drop the stray "h1-pagebreak" paragraph style from the (empty) first
#     paragraph that only carries the _GoBack bookmark.
#  2. Register two table styles in anticipation of pandoc 2.0's table
#     support: the built-in "Table Grid" style (Tabellenraster) and a
#     custom "Table" style, both based on the document's default table
#     style (Normal Table / NormaleTabelle).

$d = $word.ActiveDocument

# --- 1. Remove the h1-pagebreak paragraph style from the first paragraph ---
$first = $d.Paragraphs(1)
$first.Style = "Standard"

# --- 2. Add the "Table Grid" table style (maps to styleId Tabellenraster) ---
$tableGrid = $d.Styles.Add("Tabellenraster", 3)
$tableGrid.NameLocal = "Table Grid"
$tableGrid.BaseStyle = "NormaleTabelle"
$tableGrid.ParagraphFormat.SpaceAfter = 0

# --- 3. Add the custom "Table" table style ---
$table = $d.Styles.Add("Table", 3)
$table.NameLocal = "Table"
$table.BaseStyle = "NormaleTabelle"
$table.Priority = 99
$table.ParagraphFormat.SpaceAfter = 0
